# Reran questionnaire analysis with corrected questionnaire data.
# Updates the normality, equal_var, mixed_anova and pairwise_tests sheets
# of the MDBF interaction-questionnaire statistics workbook with the
# freshly recomputed numbers.

$wb = $excel.ActiveWorkbook

function Set-Num($ws, $cell, $val) {
    $ws.Range($cell).Value = $val
}

function Set-Txt($ws, $cell, $val) {
    # Force a genuine text cell (these BF10 entries are stored as strings,
    # e.g. "1.802e+05", not numbers) without leaving a stray number-format
    # style behind on the cell.
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# normality sheet — Shapiro-Wilk W / p-value columns (D, E)
# ---------------------------------------------------------------------
$wsNorm = $wb.Worksheets.Item("normality")

Set-Num $wsNorm "D3"  0.9314
Set-Num $wsNorm "E3"  0.3953
Set-Num $wsNorm "D4"  0.949
Set-Num $wsNorm "E4"  0.5839
Set-Num $wsNorm "D5"  0.9624
Set-Num $wsNorm "E5"  0.8173
Set-Num $wsNorm "D6"  0.9
Set-Num $wsNorm "E6"  0.1339
Set-Num $wsNorm "D7"  0.9167999999999999
Set-Num $wsNorm "E7"  0.2606
Set-Num $wsNorm "D8"  0.9529
Set-Num $wsNorm "E8"  0.6433
Set-Num $wsNorm "D9"  0.8181
Set-Num $wsNorm "E9"  0.0152
Set-Num $wsNorm "D10" 0.9473
Set-Num $wsNorm "E10" 0.5587
Set-Num $wsNorm "D11" 0.9631
Set-Num $wsNorm "E11" 0.8267
Set-Num $wsNorm "D12" 0.8824
Set-Num $wsNorm "E12" 0.0769
Set-Num $wsNorm "D13" 0.9609
Set-Num $wsNorm "E13" 0.7959000000000001
Set-Num $wsNorm "D14" 0.9126
Set-Num $wsNorm "E14" 0.1987

# ---------------------------------------------------------------------
# equal_var sheet — Levene W / p-value columns (D, E)
# ---------------------------------------------------------------------
$wsEq = $wb.Worksheets.Item("equal_var")

Set-Num $wsEq "D3" 0.2057
Set-Num $wsEq "E3" 0.6544
Set-Num $wsEq "D4" 2.9337
Set-Num $wsEq "E4" 0.1002
Set-Num $wsEq "D5" 0.6622
Set-Num $wsEq "E5" 0.4241
Set-Num $wsEq "D6" 0.1763
Set-Num $wsEq "E6" 0.6785
Set-Num $wsEq "D7" 1.8029
Set-Num $wsEq "E7" 0.1925
Set-Num $wsEq "D8" 0.2494
Set-Num $wsEq "E8" 0.6222

# ---------------------------------------------------------------------
# mixed_anova sheet — SS, MS, F, p-unc, np2 columns (D, G, H, I, J)
# ---------------------------------------------------------------------
$wsAnova = $wb.Worksheets.Item("mixed_anova")

Set-Num $wsAnova "D3"  57.8954
Set-Num $wsAnova "G3"  57.8954
Set-Num $wsAnova "H3"  1.2653
Set-Num $wsAnova "I3"  0.2723
Set-Num $wsAnova "J3"  0.0521

Set-Num $wsAnova "H4"  9.6745
Set-Num $wsAnova "I4"  0.0049
Set-Num $wsAnova "J4"  0.2961

Set-Num $wsAnova "D5"  6.3184
Set-Num $wsAnova "G5"  6.3184
Set-Num $wsAnova "H5"  0.4658
Set-Num $wsAnova "I5"  0.5017
Set-Num $wsAnova "J5"  0.0199

Set-Num $wsAnova "D6"  35.0678
Set-Num $wsAnova "G6"  35.0678
Set-Num $wsAnova "H6"  0.6408
Set-Num $wsAnova "I6"  0.4316
Set-Num $wsAnova "J6"  0.0271

Set-Num $wsAnova "H7"  58.1726
Set-Num $wsAnova "J7"  0.7167

Set-Num $wsAnova "D8"  32.6955
Set-Num $wsAnova "G8"  32.6955
Set-Num $wsAnova "H8"  1.3542
Set-Num $wsAnova "I8"  0.2565
Set-Num $wsAnova "J8"  0.0556

Set-Num $wsAnova "D9"  33.8713
Set-Num $wsAnova "G9"  33.8713
Set-Num $wsAnova "H9"  0.7938
Set-Num $wsAnova "I9"  0.3822
Set-Num $wsAnova "J9"  0.0334

Set-Num $wsAnova "H10" 69.539
Set-Num $wsAnova "J10" 0.7514999999999999

Set-Num $wsAnova "D11" 39.0646
Set-Num $wsAnova "G11" 39.0646
Set-Num $wsAnova "H11" 1.3182
Set-Num $wsAnova "I11" 0.2627
Set-Num $wsAnova "J11" 0.0542

# ---------------------------------------------------------------------
# pairwise_tests sheet — T, dof, p-unc, BF10, hedges columns (I, J, L, M, N)
# ---------------------------------------------------------------------
$wsPw = $wb.Worksheets.Item("pairwise_tests")

Set-Num $wsPw "I4"  1.1317
Set-Num $wsPw "J4"  22.9078
Set-Num $wsPw "L4"  0.2695
Set-Txt $wsPw "M4"  "0.586"
Set-Num $wsPw "N4"  0.4355

Set-Num $wsPw "I5"  1.282
Set-Num $wsPw "J5"  22.2886
Set-Num $wsPw "L5"  0.213
Set-Txt $wsPw "M5"  "0.667"
Set-Num $wsPw "N5"  0.4983

Set-Num $wsPw "I6"  0.6892
Set-Num $wsPw "J6"  20.0744
Set-Num $wsPw "L6"  0.4986
Set-Txt $wsPw "M6"  "0.438"
Set-Num $wsPw "N6"  0.2618

Set-Num $wsPw "I8"  0.8073
Set-Num $wsPw "J8"  22.6763
Set-Num $wsPw "L8"  0.4279
Set-Txt $wsPw "M8"  "0.467"
Set-Num $wsPw "N8"  0.3099

Set-Num $wsPw "I9"  1.1896
Set-Num $wsPw "J9"  21.9495
Set-Num $wsPw "L9"  0.2469
Set-Txt $wsPw "M9"  "0.615"
Set-Num $wsPw "N9"  0.4548

Set-Num $wsPw "I10" 0.0265
Set-Num $wsPw "J10" 22.9998
Set-Num $wsPw "L10" 0.9791
Set-Txt $wsPw "M10" "0.368"
Set-Num $wsPw "N10" 0.0102

Set-Num $wsPw "I12" 0.8998
Set-Num $wsPw "J12" 22.4616
Set-Num $wsPw "L12" 0.3778
Set-Txt $wsPw "M12" "0.495"
Set-Num $wsPw "N12" 0.3449

Set-Num $wsPw "I13" 1.1317
Set-Num $wsPw "J13" 21.9648
Set-Num $wsPw "L13" 0.27
Set-Txt $wsPw "M13" "0.586"
Set-Num $wsPw "N13" 0.4327

Set-Num $wsPw "I14" -0.081
Set-Num $wsPw "J14" 22.4056
Set-Num $wsPw "L14" 0.9361
Set-Txt $wsPw "M14" "0.369"
Set-Num $wsPw "N14" -0.0315
